$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 257
$ws1.Range("F5").Value = 6565
$ws1.Range("F6").Value = 5347
$ws1.Range("F9").Value = 6
$ws1.Range("F10").Value = 64
$ws1.Range("F11").Value = 229
$ws1.Range("F12").Value = 40

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 257
$ws4.Range("F5").Value = 6565
$ws4.Range("F6").Value = 5347
$ws4.Range("F9").Value = 6
$ws4.Range("F10").Value = 64
$ws4.Range("F11").Value = 229
$ws4.Range("F14").Value = 40
